$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 3: split out a large-donation entry from the general total
$ws.Range("A3").Value = 45232
$ws.Range("A3").NumberFormat = $ws.Range("A2").NumberFormat
$ws.Range("B3").Value = "General donations"
$ws.Range("C3").Value = 3000
